$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 1-3: Coca/Agua/Pan with real numeric quantities, new font style ---
$ws.Range("A1").Value = "Coca"
$ws.Range("B1").Value = 35
$ws.Range("C1").Value = 40

$ws.Range("A2").Value = "Agua"
$ws.Range("B2").Value = 15
$ws.Range("C2").Value = 50

$ws.Range("A3").Value = "Pan"
$ws.Range("B3").Value = 20
$ws.Range("C3").Value = 60

# Give rows 1-3 a distinct (blank/default) font so a new style is created
$ws.Range("A1:C3").Font.Name = ""

# --- Rows 4-6: Gas / Sprite / Sprite, quantities kept as text values ---
$ws.Range("A4").Value = "Gas"
$ws.Cells.Item(4, 2).Value = "'50"
$ws.Cells.Item(4, 3).Value = "'100"

$ws.Range("A5").Value = "Sprite"
$ws.Cells.Item(5, 2).Value = "'40"
$ws.Cells.Item(5, 3).Value = "'60"

$ws.Range("A6").Value = "Sprite"
$ws.Cells.Item(6, 2).Value = "'30"
$ws.Cells.Item(6, 3).Value = "'60"

# Drop the quote-prefix formatting so B4:C6 stay plain (unstyled) text cells
$ws.Range("B4:C6").ClearFormats()
